# Generate Report for Archive
#
# The localization status report moves from "Ready for handoff" to
# "In Translation", and the zh-cn / de-de detail sheets record the
# name of the handoff that produced this report ("TestHandoff1") in
# the "Lastest Handoff Name" column.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Status columns for zh-cn (E2) and de-de (F2)
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn detail sheet: Status (C2) and Lastest Handoff Name (I2)
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("I2").Value = "TestHandoff1"

# de-de detail sheet: Status (C2) and Lastest Handoff Name (I2)
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("I2").Value = "TestHandoff1"

# Column widths on the Status columns shrink slightly to fit the new,
# shorter "In Translation" text (previously sized for "Ready for handoff").
$wsOverview.Columns.Item(5).ColumnWidth = 12.42
$wsOverview.Columns.Item(6).ColumnWidth = 12.42
$wsZhCn.Columns.Item(3).ColumnWidth = 12.42
$wsDeDe.Columns.Item(3).ColumnWidth = 12.42
